$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.561187267303467
$ws.Range("B1").Value = 2.259199142456055
$ws.Range("C1").Value = 3.374518156051636
$ws.Range("D1").Value = 2.442772626876831
$ws.Range("E1").Value = 0.7161096930503845
